$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1 (cg_AAA): insert a new coverpoint "cp_rxpkt_len_trans"
# (ID 6) right before the existing "cc_pktlen_port_en" coverpoint,
# which becomes ID 7 and shifts two rows down.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Make room: push everything from row 17 onward down by two rows.
$ws1.Rows("17:18").Insert()

# New coverpoint block (rows 17-18): transition bins on rxpkt_len.
$ws1.Range("B17").Value = 6
$ws1.Range("C17").Value = "cp_rxpkt_len_trans"
$ws1.Range("D17").Value = "pkt_en"
$ws1.Range("G17").Value = "MIN_MAX"
$ws1.Range("H17").Value = "64 => 1518"
$ws1.Range("I17").Value = "bm.pkt"

$ws1.Range("G18").Value = "MAX_MIN"
$ws1.Range("H18").Value = "1518 => 64"
$ws1.Range("I18").Value = "bm.pkt"

# The old row 17 ("cc_pktlen_port_en") is now row 19 - renumber it
# from ID 6 to ID 7 (everything else already carried down by Insert).
$ws1.Range("B19").Value = 7

$ws1.Range("A2:XFD2").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet 2 (cg_BBB): replace the old coverage content with a single
# "cp_chipmode" coverpoint (MODE0..MODE3 + illegal_bins).
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2:K21").ClearContents()

$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = "cp_chipmode"
$ws2.Range("D2").Value = "mode"
$ws2.Range("G2").Value = "MODE0"
$ws2.Range("H2").Value = "2'b00"
$ws2.Range("I2").Value = " "
$ws2.Range("J2").Value = " "
$ws2.Range("K2").Value = " "

$ws2.Range("G3").Value = "MODE1"
$ws2.Range("H3").Value = "2'b01"

$ws2.Range("G4").Value = "MODE2"
$ws2.Range("H4").Value = "2'b10"

$ws2.Range("F5").Value = "illegal_bins"
$ws2.Range("G5").Value = "MODE3"
$ws2.Range("H5").Value = "2'b11"

$ws2.Range("A9").Value = "$"

$ws2.Range("C21").Value = " "
$ws2.Range("A23").Value = " "

$ws2.Columns("G").ColumnWidth = 18.36328125

$ws2.Range("A21").Select() | Out-Null
